$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 191  # ALC H9: 189.5 -> 191
$ws.Cells.Item(9, 10).Value = 191.25  # ALC J9: 189.4 -> 191.25
$ws.Cells.Item(9, 12).Value = 191.25  # ALC L9: 189.4 -> 191.25
$ws.Cells.Item(9, 14).Value = -529.25  # ALC N9: -527.4 -> -529.25

$ws.Cells.Item(12, 8).Value = 899.6667  # ALC H12: 849.75 -> 899.6667
$ws.Cells.Item(12, 9).Value = 849.5  # ALC I12: 799.6667 -> 849.5
$ws.Cells.Item(12, 11).Value = 849.5  # ALC K12: 799.6667 -> 849.5
$ws.Cells.Item(12, 13).Value = -679.5  # ALC M12: -629.6667 -> -679.5

$ws.Cells.Item(40, 8).Value = 3381.8235  # ALC H40: 3522.1 -> 3381.8235
$ws.Cells.Item(40, 9).Value = 2494.6667  # ALC I40: 3000 -> 2494.6667
$ws.Cells.Item(40, 10).Value = 3571.9285  # ALC J40: 3549.5789 -> 3571.9285
$ws.Cells.Item(40, 11).Value = 2494.6667  # ALC K40: 3000 -> 2494.6667
$ws.Cells.Item(40, 12).Value = 3571.9285  # ALC L40: 3549.5789 -> 3571.9285
$ws.Cells.Item(40, 13).Value = -2319.6667  # ALC M40: -2825 -> -2319.6667
$ws.Cells.Item(40, 14).Value = -3921.9285  # ALC N40: -3899.5789 -> -3921.9285

$ws.Cells.Item(43, 9).Value = 999.5  # ALC I43: 999.3333 -> 999.5
$ws.Cells.Item(43, 10).Value = 999  # ALC J43: 0 -> 999
$ws.Cells.Item(43, 11).Value = 999.5  # ALC K43: 999.3333 -> 999.5
$ws.Cells.Item(43, 12).Value = 999  # ALC L43: 0 -> 999
$ws.Cells.Item(43, 13).Value = -930.5  # ALC M43: -930.3333 -> -930.5
$ws.Cells.Item(43, 14).Value = -1137  # ALC N43: None -> -1137

$ws.Cells.Item(58, 8).Value = 3396.6667  # ALC H58: 3197 -> 3396.6667
$ws.Cells.Item(58, 10).Value = 3396.6667  # ALC J58: 3197 -> 3396.6667
$ws.Cells.Item(58, 12).Value = 10190.0001  # ALC L58: 9591 -> 10190.0001
$ws.Cells.Item(58, 14).Value = -10490.0001  # ALC N58: -9891 -> -10490.0001

$ws.Cells.Item(98, 8).Value = 1953.625  # ALC H98: 1995.25 -> 1953.625
$ws.Cells.Item(98, 9).Value = 670.3333  # ALC I98: 743.625 -> 670.3333
$ws.Cells.Item(98, 10).Value = 3603.5715  # ALC J98: 3246.875 -> 3603.5715
$ws.Cells.Item(98, 11).Value = 670.3333  # ALC K98: 743.625 -> 670.3333
$ws.Cells.Item(98, 12).Value = 3603.5715  # ALC L98: 3246.875 -> 3603.5715
$ws.Cells.Item(98, 13).Value = 827.6667  # ALC M98: 754.375 -> 827.6667
$ws.Cells.Item(98, 14).Value = -6599.5715  # ALC N98: -6242.875 -> -6599.5715

$ws.Cells.Item(111, 8).Value = 1500  # ALC H111: 400 -> 1500
$ws.Cells.Item(111, 9).Value = 1500  # ALC I111: 0 -> 1500
$ws.Cells.Item(111, 10).Value = 0  # ALC J111: 400 -> 0
$ws.Cells.Item(111, 11).Value = 4500  # ALC K111: 0 -> 4500
$ws.Cells.Item(111, 12).Value = 0  # ALC L111: 1200 -> 0
$ws.Cells.Item(111, 13).Value = -1433  # ALC M111: None -> -1433
$ws.Cells.Item(111, 14).ClearContents()  # ALC N111: -7334 -> (removed)

$ws.Cells.Item(122, 8).Value = 1953.625  # ALC H122: 1995.25 -> 1953.625
$ws.Cells.Item(122, 9).Value = 670.3333  # ALC I122: 743.625 -> 670.3333
$ws.Cells.Item(122, 10).Value = 3603.5715  # ALC J122: 3246.875 -> 3603.5715
$ws.Cells.Item(122, 11).Value = 2010.9999  # ALC K122: 2230.875 -> 2010.9999
$ws.Cells.Item(122, 12).Value = 10810.7145  # ALC L122: 9740.625 -> 10810.7145
$ws.Cells.Item(122, 13).Value = 439.0001  # ALC M122: 219.125 -> 439.0001
$ws.Cells.Item(122, 14).Value = -15710.7145  # ALC N122: -14640.625 -> -15710.7145

$ws.Cells.Item(132, 8).Value = 38465972  # ALC H132: 38466210 -> 38465972
$ws.Cells.Item(132, 9).Value = 50004710  # ALC I132: 50005020 -> 50004710
$ws.Cells.Item(132, 11).Value = 150014130  # ALC K132: 150015060 -> 150014130
$ws.Cells.Item(132, 13).Value = -150011600  # ALC M132: -150012530 -> -150011600

$ws.Cells.Item(138, 8).Value = 3530.5881  # ALC H138: 3734.1667 -> 3530.5881
$ws.Cells.Item(138, 10).Value = 3601.4  # ALC J138: 3881.1 -> 3601.4
$ws.Cells.Item(138, 12).Value = 10804.2  # ALC L138: 11643.3 -> 10804.2
$ws.Cells.Item(138, 14).Value = -21084.2  # ALC N138: -21923.3 -> -21084.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6927.875  # ARM H32: 6940.35 -> 6927.875
$ws.Cells.Item(32, 9).Value = 6927.875  # ARM I32: 6940.35 -> 6927.875
$ws.Cells.Item(32, 11).Value = 6927.875  # ARM K32: 6940.35 -> 6927.875
$ws.Cells.Item(32, 13).Value = -6640.875  # ARM M32: -6653.35 -> -6640.875

$ws.Cells.Item(61, 8).Value = 806  # ARM H61: 737.3333 -> 806
$ws.Cells.Item(61, 9).Value = 806  # ARM I61: 737.3333 -> 806
$ws.Cells.Item(61, 11).Value = 806  # ARM K61: 737.3333 -> 806
$ws.Cells.Item(61, 13).Value = -594  # ARM M61: -525.3333 -> -594

$ws.Cells.Item(74, 8).Value = 994.5  # ARM H74: 799.6667 -> 994.5
$ws.Cells.Item(74, 9).Value = 994.5  # ARM I74: 799.6667 -> 994.5
$ws.Cells.Item(74, 11).Value = 994.5  # ARM K74: 799.6667 -> 994.5
$ws.Cells.Item(74, 13).Value = -120.5  # ARM M74: 74.33330000000001 -> -120.5

$ws.Cells.Item(77, 8).Value = 994.5  # ARM H77: 799.6667 -> 994.5
$ws.Cells.Item(77, 9).Value = 994.5  # ARM I77: 799.6667 -> 994.5
$ws.Cells.Item(77, 11).Value = 4972.5  # ARM K77: 3998.3335 -> 4972.5
$ws.Cells.Item(77, 13).Value = -604.5  # ARM M77: 369.6665000000003 -> -604.5

$ws.Cells.Item(80, 8).Value = 40000  # ARM H80: 38333.332 -> 40000
$ws.Cells.Item(80, 9).Value = 0  # ARM I80: 30000 -> 0
$ws.Cells.Item(80, 11).Value = 0  # ARM K80: 30000 -> 0
$ws.Cells.Item(80, 13).ClearContents()  # ARM M80: -29002 -> (removed)

$ws.Cells.Item(83, 8).Value = 40000  # ARM H83: 38333.332 -> 40000
$ws.Cells.Item(83, 9).Value = 0  # ARM I83: 30000 -> 0
$ws.Cells.Item(83, 11).Value = 0  # ARM K83: 90000 -> 0
$ws.Cells.Item(83, 13).ClearContents()  # ARM M83: -85008 -> (removed)

$ws.Cells.Item(136, 8).Value = 806  # ARM H136: 737.3333 -> 806
$ws.Cells.Item(136, 9).Value = 806  # ARM I136: 737.3333 -> 806
$ws.Cells.Item(136, 11).Value = 2418  # ARM K136: 2211.9999 -> 2418
$ws.Cells.Item(136, 13).Value = 132  # ARM M136: 338.0001000000002 -> 132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 1094.5  # BSM H37: 1159.3334 -> 1094.5
$ws.Cells.Item(37, 9).Value = 1094.5  # BSM I37: 1159.3334 -> 1094.5
$ws.Cells.Item(37, 11).Value = 1094.5  # BSM K37: 1159.3334 -> 1094.5
$ws.Cells.Item(37, 13).Value = -957.5  # BSM M37: -1022.3334 -> -957.5

$ws.Cells.Item(134, 8).Value = 5894  # BSM H134: 6543.5625 -> 5894
$ws.Cells.Item(134, 9).Value = 5894  # BSM I134: 6559.933 -> 5894
$ws.Cells.Item(134, 10).Value = 0  # BSM J134: 6298 -> 0
$ws.Cells.Item(134, 11).Value = 17682  # BSM K134: 19679.799 -> 17682
$ws.Cells.Item(134, 12).Value = 0  # BSM L134: 18894 -> 0
$ws.Cells.Item(134, 13).Value = -15147  # BSM M134: -17144.799 -> -15147
$ws.Cells.Item(134, 14).ClearContents()  # BSM N134: -23964 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 14284.667  # CRP H41: 12233.818 -> 14284.667
$ws.Cells.Item(41, 9).Value = 8712.4  # CRP I41: 7081.7144 -> 8712.4
$ws.Cells.Item(41, 11).Value = 8712.4  # CRP K41: 7081.7144 -> 8712.4
$ws.Cells.Item(41, 13).Value = -8284.4  # CRP M41: -6653.7144 -> -8284.4

$ws.Cells.Item(59, 8).Value = 27077.908  # CRP H59: 27539.777 -> 27077.908
$ws.Cells.Item(59, 9).Value = 22551  # CRP I59: 23809.666 -> 22551
$ws.Cells.Item(59, 11).Value = 22551  # CRP K59: 23809.666 -> 22551
$ws.Cells.Item(59, 13).Value = -21406  # CRP M59: -22664.666 -> -21406

$ws.Cells.Item(60, 8).Value = 9935.875  # CRP H60: 8410.875 -> 9935.875
$ws.Cells.Item(60, 9).Value = 4914.5  # CRP I60: 2881.1667 -> 4914.5
$ws.Cells.Item(60, 11).Value = 4914.5  # CRP K60: 2881.1667 -> 4914.5
$ws.Cells.Item(60, 13).Value = -4403.5  # CRP M60: -2370.1667 -> -4403.5

$ws.Cells.Item(134, 8).Value = 4598.5713  # CRP H134: 4798.4614 -> 4598.5713
$ws.Cells.Item(134, 9).Value = 4217.4  # CRP I134: 4463.778 -> 4217.4
$ws.Cells.Item(134, 11).Value = 12652.2  # CRP K134: 13391.334 -> 12652.2
$ws.Cells.Item(134, 13).Value = -10117.2  # CRP M134: -10856.334 -> -10117.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 2999.3333  # CUL H80: 2501 -> 2999.3333
$ws.Cells.Item(80, 9).Value = 0  # CUL I80: 2002 -> 0
$ws.Cells.Item(80, 10).Value = 2999.3333  # CUL J80: 3000 -> 2999.3333
$ws.Cells.Item(80, 11).Value = 0  # CUL K80: 6006 -> 0
$ws.Cells.Item(80, 12).Value = 8997.999899999999  # CUL L80: 9000 -> 8997.999899999999
$ws.Cells.Item(80, 13).ClearContents()  # CUL M80: -5070 -> (removed)
$ws.Cells.Item(80, 14).Value = -10869.9999  # CUL N80: -10872 -> -10869.9999

$ws.Cells.Item(83, 8).Value = 2999.3333  # CUL H83: 2501 -> 2999.3333
$ws.Cells.Item(83, 9).Value = 0  # CUL I83: 2002 -> 0
$ws.Cells.Item(83, 10).Value = 2999.3333  # CUL J83: 3000 -> 2999.3333
$ws.Cells.Item(83, 11).Value = 0  # CUL K83: 18018 -> 0
$ws.Cells.Item(83, 12).Value = 26993.9997  # CUL L83: 27000 -> 26993.9997
$ws.Cells.Item(83, 13).ClearContents()  # CUL M83: -13338 -> (removed)
$ws.Cells.Item(83, 14).Value = -36353.9997  # CUL N83: -36360 -> -36353.9997

$ws.Cells.Item(131, 8).Value = 2201.5  # CUL H131: 2269.111 -> 2201.5
$ws.Cells.Item(131, 9).Value = 1701.6666  # CUL I131: 1715.25 -> 1701.6666
$ws.Cells.Item(131, 11).Value = 5104.9998  # CUL K131: 5145.75 -> 5104.9998
$ws.Cells.Item(131, 13).Value = -64.9997999999996  # CUL M131: -105.75 -> -64.9997999999996

$ws.Cells.Item(137, 8).Value = 10600  # CUL H137: 7498.6665 -> 10600
$ws.Cells.Item(137, 10).Value = 13466.667  # CUL J137: 8598.4 -> 13466.667
$ws.Cells.Item(137, 12).Value = 40400.001  # CUL L137: 25795.2 -> 40400.001
$ws.Cells.Item(137, 14).Value = -50600.001  # CUL N137: -35995.2 -> -50600.001

$ws.Cells.Item(138, 8).Value = 1774.5  # CUL H138: 1613.8572 -> 1774.5
$ws.Cells.Item(138, 9).Value = 647  # CUL I138: 648.5 -> 647
$ws.Cells.Item(138, 11).Value = 1941  # CUL K138: 1945.5 -> 1941
$ws.Cells.Item(138, 13).Value = 3199  # CUL M138: 3194.5 -> 3199

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 40013.555  # GSM H122: 37542.15 -> 40013.555
$ws.Cells.Item(122, 9).Value = 46228.848  # GSM I122: 42969.645 -> 46228.848
$ws.Cells.Item(122, 10).Value = 23853.8  # GSM J122: 24878 -> 23853.8
$ws.Cells.Item(122, 11).Value = 138686.544  # GSM K122: 128908.935 -> 138686.544
$ws.Cells.Item(122, 12).Value = 71561.39999999999  # GSM L122: 74634 -> 71561.39999999999
$ws.Cells.Item(122, 13).Value = -136236.544  # GSM M122: -126458.935 -> -136236.544
$ws.Cells.Item(122, 14).Value = -76461.39999999999  # GSM N122: -79534 -> -76461.39999999999

$ws.Cells.Item(132, 8).Value = 1749.6  # GSM H132: 1987 -> 1749.6
$ws.Cells.Item(132, 9).Value = 1749.6  # GSM I132: 1987 -> 1749.6
$ws.Cells.Item(132, 11).Value = 5248.799999999999  # GSM K132: 5961 -> 5248.799999999999
$ws.Cells.Item(132, 13).Value = -2718.799999999999  # GSM M132: -3431 -> -2718.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3340.5  # LTW H40: 3426.7058 -> 3340.5
$ws.Cells.Item(40, 9).Value = 3408.6  # LTW I40: 3518.1428 -> 3408.6
$ws.Cells.Item(40, 11).Value = 3408.6  # LTW K40: 3518.1428 -> 3408.6
$ws.Cells.Item(40, 13).Value = -3272.6  # LTW M40: -3382.1428 -> -3272.6

$ws.Cells.Item(58, 8).Value = 2848.75  # LTW H58: 2997.5 -> 2848.75
$ws.Cells.Item(58, 9).Value = 2848.75  # LTW I58: 2997.5 -> 2848.75
$ws.Cells.Item(58, 11).Value = 2848.75  # LTW K58: 2997.5 -> 2848.75
$ws.Cells.Item(58, 13).Value = -2588.75  # LTW M58: -2737.5 -> -2588.75

$ws.Cells.Item(82, 8).Value = 2537.75  # LTW H82: 2230.2 -> 2537.75
$ws.Cells.Item(82, 9).Value = 2537.75  # LTW I82: 2230.2 -> 2537.75
$ws.Cells.Item(82, 11).Value = 2537.75  # LTW K82: 2230.2 -> 2537.75
$ws.Cells.Item(82, 13).Value = -2176.75  # LTW M82: -1869.2 -> -2176.75

$ws.Cells.Item(85, 8).Value = 2537.75  # LTW H85: 2230.2 -> 2537.75
$ws.Cells.Item(85, 9).Value = 2537.75  # LTW I85: 2230.2 -> 2537.75
$ws.Cells.Item(85, 11).Value = 2537.75  # LTW K85: 2230.2 -> 2537.75
$ws.Cells.Item(85, 13).Value = -1289.75  # LTW M85: -982.1999999999998 -> -1289.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 0  # WVR H81: 553.8 -> 0
$ws.Cells.Item(81, 9).Value = 0  # WVR I81: 553.8 -> 0
$ws.Cells.Item(81, 11).Value = 0  # WVR K81: 1107.6 -> 0
$ws.Cells.Item(81, 13).ClearContents()  # WVR M81: -46.59999999999991 -> (removed)

$ws.Cells.Item(84, 8).Value = 0  # WVR H84: 553.8 -> 0
$ws.Cells.Item(84, 9).Value = 0  # WVR I84: 553.8 -> 0
$ws.Cells.Item(84, 11).Value = 0  # WVR K84: 5538 -> 0
$ws.Cells.Item(84, 13).ClearContents()  # WVR M84: -234 -> (removed)

$ws.Cells.Item(104, 8).Value = 27124.5  # WVR H104: 29499.666 -> 27124.5
$ws.Cells.Item(104, 10).Value = 27124.5  # WVR J104: 29499.666 -> 27124.5
$ws.Cells.Item(104, 12).Value = 27124.5  # WVR L104: 29499.666 -> 27124.5
$ws.Cells.Item(104, 14).Value = -34112.5  # WVR N104: -36487.666 -> -34112.5

$ws.Cells.Item(136, 8).Value = 3308  # WVR H136: 1742.0834 -> 3308
$ws.Cells.Item(136, 9).Value = 3308  # WVR I136: 1742.0834 -> 3308
$ws.Cells.Item(136, 11).Value = 9924  # WVR K136: 5226.2502 -> 9924
$ws.Cells.Item(136, 13).Value = -7374  # WVR M136: -2676.2502 -> -7374
